$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension will auto-expand when row 192 is written.

# Row 95
$ws.Range("D95").Value = 44966
$ws.Range("H95").Value = 'Sin especificar'
$ws.Range("J95").Value = 100
$ws.Range("K95").Value = 8500
$ws.Range("L95").Value = 9000
$ws.Range("M95").Value = 8750
$ws.Range("N95").Value = '$/caja 50 unidades'
$ws.Range("O95").Value = 'Región de O''Higgins'
$ws.Range("P95").Value = 175
$ws.Range("Q95").Value = 50

# Row 96
$ws.Range("D96").Value = 44932
$ws.Range("H96").Value = 'Sin especificar'
$ws.Range("J96").Value = 450
$ws.Range("K96").Value = 10000
$ws.Range("L96").Value = 11000
$ws.Range("M96").Value = 10556
$ws.Range("N96").Value = '$/caja 60 unidades'
$ws.Range("O96").Value = 'Región de Arica y Parinacota'
$ws.Range("P96").Value = 176
$ws.Range("Q96").Value = 60

# Row 97
$ws.Range("D97").Value = 44336
$ws.Range("H97").Value = 'Sin especificar'
$ws.Range("J97").Value = 100
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 12000
$ws.Range("M97").Value = 11000
$ws.Range("N97").Value = '$/caja 50 unidades'
$ws.Range("O97").Value = 'Región de Arica y Parinacota'
$ws.Range("P97").Value = 220
$ws.Range("Q97").Value = 50

# Row 98
$ws.Range("D98").Value = 44551
$ws.Range("H98").Value = 'Sin especificar'
$ws.Range("J98").Value = 100
$ws.Range("K98").Value = 9000
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 9500
$ws.Range("N98").Value = '$/caja 50 unidades'
$ws.Range("O98").Value = 'Región de O''Higgins'
$ws.Range("P98").Value = 190
$ws.Range("Q98").Value = 50

# Row 99
$ws.Range("D99").Value = 44203
$ws.Range("H99").Value = 'Sin especificar'
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 10000
$ws.Range("L99").Value = 11000
$ws.Range("M99").Value = 10500
$ws.Range("N99").Value = '$/caja 60 unidades'
$ws.Range("O99").Value = 'Región de O''Higgins'
$ws.Range("P99").Value = 175
$ws.Range("Q99").Value = 60

# Row 100
$ws.Range("D100").Value = 44258
$ws.Range("H100").Value = 'Sin especificar'
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 11000
$ws.Range("M100").Value = 10500
$ws.Range("N100").Value = '$/caja 60 unidades'
$ws.Range("O100").Value = 'Región de O''Higgins'
$ws.Range("P100").Value = 175
$ws.Range("Q100").Value = 60

# Row 101
$ws.Range("D101").Value = 44673
$ws.Range("H101").Value = 'Huracán'
$ws.Range("J101").Value = 350
$ws.Range("K101").Value = 9000
$ws.Range("L101").Value = 9500
$ws.Range("M101").Value = 9286
$ws.Range("N101").Value = '$/caja 60 unidades'
$ws.Range("O101").Value = 'Región de Arica y Parinacota'
$ws.Range("P101").Value = 155
$ws.Range("Q101").Value = 60

# Row 102
$ws.Range("D102").Value = 44665
$ws.Range("H102").Value = 'Sin especificar'
$ws.Range("J102").Value = 200
$ws.Range("K102").Value = 11000
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = 11500
$ws.Range("N102").Value = '$/caja 50 unidades'
$ws.Range("O102").Value = 'Región de Arica y Parinacota'
$ws.Range("P102").Value = 230
$ws.Range("Q102").Value = 50

# Row 103
$ws.Range("D103").Value = 44608
$ws.Range("H103").Value = 'Sin especificar'
$ws.Range("J103").Value = 100
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = 9500
$ws.Range("N103").Value = '$/caja 50 unidades'
$ws.Range("O103").Value = 'Región de O''Higgins'
$ws.Range("P103").Value = 190
$ws.Range("Q103").Value = 50

# Row 104
$ws.Range("D104").Value = 44721
$ws.Range("H104").Value = 'Huracán'
$ws.Range("J104").Value = 180
$ws.Range("K104").Value = 7500
$ws.Range("L104").Value = 8000
$ws.Range("M104").Value = 7778
$ws.Range("N104").Value = '$/caja 60 unidades'
$ws.Range("O104").Value = 'Región de Arica y Parinacota'
$ws.Range("P104").Value = 130
$ws.Range("Q104").Value = 60

# Row 105
$ws.Range("D105").Value = 44211
$ws.Range("H105").Value = 'Sin especificar'
$ws.Range("J105").Value = 200
$ws.Range("K105").Value = 11000
$ws.Range("L105").Value = 12000
$ws.Range("M105").Value = 11500
$ws.Range("N105").Value = '$/caja 60 unidades'
$ws.Range("O105").Value = 'Región de O''Higgins'
$ws.Range("P105").Value = 192
$ws.Range("Q105").Value = 60

# Row 106
$ws.Range("D106").Value = 44546
$ws.Range("H106").Value = 'Sin especificar'
$ws.Range("J106").Value = 250
$ws.Range("K106").Value = 8000
$ws.Range("L106").Value = 9000
$ws.Range("M106").Value = 8400
$ws.Range("N106").Value = '$/caja 60 unidades'
$ws.Range("O106").Value = 'Región Metropolitana'
$ws.Range("P106").Value = 140
$ws.Range("Q106").Value = 60

# Row 107
$ws.Range("D107").Value = 44755
$ws.Range("H107").Value = 'Sin especificar'
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 8000
$ws.Range("L107").Value = 8500
$ws.Range("M107").Value = 8250
$ws.Range("N107").Value = '$/caja 50 unidades'
$ws.Range("O107").Value = 'Región de Arica y Parinacota'
$ws.Range("P107").Value = 165
$ws.Range("Q107").Value = 50

# Row 108
$ws.Range("D108").Value = 44580
$ws.Range("H108").Value = 'Sin especificar'
$ws.Range("J108").Value = 100
$ws.Range("K108").Value = 14000
$ws.Range("L108").Value = 15000
$ws.Range("M108").Value = 14500
$ws.Range("N108").Value = '$/caja 50 unidades'
$ws.Range("O108").Value = 'Región de O''Higgins'
$ws.Range("P108").Value = 290
$ws.Range("Q108").Value = 50

# Row 109
$ws.Range("D109").Value = 44517
$ws.Range("H109").Value = 'Sin especificar'
$ws.Range("J109").Value = 250
$ws.Range("K109").Value = 6500
$ws.Range("L109").Value = 7000
$ws.Range("M109").Value = 6700
$ws.Range("N109").Value = '$/caja 60 unidades'
$ws.Range("O109").Value = 'Región del Maule'
$ws.Range("P109").Value = 112
$ws.Range("Q109").Value = 60

# Row 110
$ws.Range("D110").Value = 44635
$ws.Range("H110").Value = 'Sin especificar'
$ws.Range("J110").Value = 220
$ws.Range("K110").Value = 12000
$ws.Range("L110").Value = 13000
$ws.Range("M110").Value = 12545
$ws.Range("N110").Value = '$/caja 60 unidades'
$ws.Range("O110").Value = 'Región Metropolitana'
$ws.Range("P110").Value = 209
$ws.Range("Q110").Value = 60

# Row 111
$ws.Range("D111").Value = 44397
$ws.Range("H111").Value = 'Sin especificar'
$ws.Range("J111").Value = 100
$ws.Range("K111").Value = 8000
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = 8500
$ws.Range("N111").Value = '$/caja 50 unidades'
$ws.Range("O111").Value = 'Región de Arica y Parinacota'
$ws.Range("P111").Value = 170
$ws.Range("Q111").Value = 50

# Row 112
$ws.Range("D112").Value = 44831
$ws.Range("H112").Value = 'Sin especificar'
$ws.Range("J112").Value = 270
$ws.Range("K112").Value = 18000
$ws.Range("L112").Value = 19000
$ws.Range("M112").Value = 18444
$ws.Range("N112").Value = '$/caja 50 unidades'
$ws.Range("O112").Value = 'Región de Arica y Parinacota'
$ws.Range("P112").Value = 369
$ws.Range("Q112").Value = 50

# Row 113
$ws.Range("D113").Value = 44420
$ws.Range("H113").Value = 'Sin especificar'
$ws.Range("J113").Value = 100
$ws.Range("K113").Value = 9000
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 9500
$ws.Range("N113").Value = '$/caja 50 unidades'
$ws.Range("O113").Value = 'Región de Arica y Parinacota'
$ws.Range("P113").Value = 190
$ws.Range("Q113").Value = 50

# Row 114
$ws.Range("D114").Value = 44238
$ws.Range("H114").Value = 'Sin especificar'
$ws.Range("J114").Value = 100
$ws.Range("K114").Value = 10000
$ws.Range("L114").Value = 11000
$ws.Range("M114").Value = 10500
$ws.Range("N114").Value = '$/caja 60 unidades'
$ws.Range("O114").Value = 'Región de O''Higgins'
$ws.Range("P114").Value = 175
$ws.Range("Q114").Value = 60

# Row 115
$ws.Range("D115").Value = 44600
$ws.Range("H115").Value = 'Sin especificar'
$ws.Range("J115").Value = 180
$ws.Range("K115").Value = 8000
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = 8556
$ws.Range("N115").Value = '$/caja 60 unidades'
$ws.Range("O115").Value = 'Región de Arica y Parinacota'
$ws.Range("P115").Value = 143
$ws.Range("Q115").Value = 60

# Row 116
$ws.Range("D116").Value = 44714
$ws.Range("H116").Value = 'Sin especificar'
$ws.Range("J116").Value = 200
$ws.Range("K116").Value = 11000
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = 11500
$ws.Range("N116").Value = '$/caja 50 unidades'
$ws.Range("O116").Value = 'Región de Arica y Parinacota'
$ws.Range("P116").Value = 230
$ws.Range("Q116").Value = 50

# Row 117
$ws.Range("D117").Value = 44558
$ws.Range("H117").Value = 'Sin especificar'
$ws.Range("J117").Value = 250
$ws.Range("K117").Value = 6000
$ws.Range("L117").Value = 6500
$ws.Range("M117").Value = 6200
$ws.Range("N117").Value = '$/caja 60 unidades'
$ws.Range("O117").Value = 'Región de Arica y Parinacota'
$ws.Range("P117").Value = 103
$ws.Range("Q117").Value = 60

# Row 118
$ws.Range("D118").Value = 44883
$ws.Range("H118").Value = 'Sin especificar'
$ws.Range("J118").Value = 100
$ws.Range("K118").Value = 10000
$ws.Range("L118").Value = 11000
$ws.Range("M118").Value = 10500
$ws.Range("N118").Value = '$/caja 50 unidades'
$ws.Range("O118").Value = 'Región de O''Higgins'
$ws.Range("P118").Value = 210
$ws.Range("Q118").Value = 50

# Row 119
$ws.Range("D119").Value = 44957
$ws.Range("H119").Value = 'Sin especificar'
$ws.Range("J119").Value = 220
$ws.Range("K119").Value = 5500
$ws.Range("L119").Value = 6000
$ws.Range("M119").Value = 5773
$ws.Range("N119").Value = '$/caja 60 unidades'
$ws.Range("O119").Value = 'Región de O''Higgins'
$ws.Range("P119").Value = 96
$ws.Range("Q119").Value = 60

# Row 120
$ws.Range("D120").Value = 44355
$ws.Range("H120").Value = 'Sin especificar'
$ws.Range("J120").Value = 100
$ws.Range("K120").Value = 9000
$ws.Range("L120").Value = 10000
$ws.Range("M120").Value = 9500
$ws.Range("N120").Value = '$/caja 50 unidades'
$ws.Range("O120").Value = 'Región de Arica y Parinacota'
$ws.Range("P120").Value = 190
$ws.Range("Q120").Value = 50

# Row 121
$ws.Range("D121").Value = 44299
$ws.Range("H121").Value = 'Sin especificar'
$ws.Range("J121").Value = 100
$ws.Range("K121").Value = 7000
$ws.Range("L121").Value = 8000
$ws.Range("M121").Value = 7500
$ws.Range("N121").Value = '$/caja 50 unidades'
$ws.Range("O121").Value = 'Región Metropolitana'
$ws.Range("P121").Value = 150
$ws.Range("Q121").Value = 50

# Row 122
$ws.Range("D122").Value = 44915
$ws.Range("H122").Value = 'Sin especificar'
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 7000
$ws.Range("L122").Value = 8000
$ws.Range("M122").Value = 7500
$ws.Range("N122").Value = '$/caja 50 unidades'
$ws.Range("O122").Value = 'Región de O''Higgins'
$ws.Range("P122").Value = 150
$ws.Range("Q122").Value = 50

# Row 123
$ws.Range("D123").Value = 44736
$ws.Range("H123").Value = 'Sin especificar'
$ws.Range("J123").Value = 100
$ws.Range("K123").Value = 16000
$ws.Range("L123").Value = 17000
$ws.Range("M123").Value = 16500
$ws.Range("N123").Value = '$/caja 60 unidades'
$ws.Range("O123").Value = 'Región de Arica y Parinacota'
$ws.Range("P123").Value = 275
$ws.Range("Q123").Value = 60

# Row 124
$ws.Range("D124").Value = 44642
$ws.Range("H124").Value = 'Sin especificar'
$ws.Range("J124").Value = 220
$ws.Range("K124").Value = 13000
$ws.Range("L124").Value = 14000
$ws.Range("M124").Value = 13455
$ws.Range("N124").Value = '$/caja 60 unidades'
$ws.Range("O124").Value = 'Región de Arica y Parinacota'
$ws.Range("P124").Value = 224
$ws.Range("Q124").Value = 60

# Row 125
$ws.Range("D125").Value = 44649
$ws.Range("H125").Value = 'Sin especificar'
$ws.Range("J125").Value = 220
$ws.Range("K125").Value = 12000
$ws.Range("L125").Value = 13000
$ws.Range("M125").Value = 12455
$ws.Range("N125").Value = '$/caja 60 unidades'
$ws.Range("O125").Value = 'Región de Arica y Parinacota'
$ws.Range("P125").Value = 208
$ws.Range("Q125").Value = 60

# Row 126
$ws.Range("D126").Value = 44278
$ws.Range("H126").Value = 'Sin especificar'
$ws.Range("J126").Value = 100
$ws.Range("K126").Value = 8000
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = 8500
$ws.Range("N126").Value = '$/caja 60 unidades'
$ws.Range("O126").Value = 'Región de O''Higgins'
$ws.Range("P126").Value = 142
$ws.Range("Q126").Value = 60

# Row 127
$ws.Range("D127").Value = 44293
$ws.Range("H127").Value = 'Sin especificar'
$ws.Range("J127").Value = 100
$ws.Range("K127").Value = 8000
$ws.Range("L127").Value = 9000
$ws.Range("M127").Value = 8500
$ws.Range("N127").Value = '$/caja 60 unidades'
$ws.Range("O127").Value = 'Región del Maule'
$ws.Range("P127").Value = 142
$ws.Range("Q127").Value = 60

# Row 128
$ws.Range("D128").Value = 44379
$ws.Range("H128").Value = 'Sin especificar'
$ws.Range("J128").Value = 100
$ws.Range("K128").Value = 9000
$ws.Range("L128").Value = 10000
$ws.Range("M128").Value = 9500
$ws.Range("N128").Value = '$/caja 50 unidades'
$ws.Range("O128").Value = 'Región de Arica y Parinacota'
$ws.Range("P128").Value = 190
$ws.Range("Q128").Value = 50

# Row 129
$ws.Range("D129").Value = 44747
$ws.Range("H129").Value = 'Sin especificar'
$ws.Range("J129").Value = 250
$ws.Range("K129").Value = 12000
$ws.Range("L129").Value = 13000
$ws.Range("M129").Value = 12400
$ws.Range("N129").Value = '$/caja 50 unidades'
$ws.Range("O129").Value = 'Región de Arica y Parinacota'
$ws.Range("P129").Value = 248
$ws.Range("Q129").Value = 50

# Row 130
$ws.Range("D130").Value = 44447
$ws.Range("H130").Value = 'Sin especificar'
$ws.Range("J130").Value = 100
$ws.Range("K130").Value = 15000
$ws.Range("L130").Value = 16000
$ws.Range("M130").Value = 15500
$ws.Range("N130").Value = '$/caja 50 unidades'
$ws.Range("O130").Value = 'Región de Arica y Parinacota'
$ws.Range("P130").Value = 310
$ws.Range("Q130").Value = 50

# Row 131
$ws.Range("D131").Value = 44434
$ws.Range("H131").Value = 'Sin especificar'
$ws.Range("J131").Value = 100
$ws.Range("K131").Value = 12000
$ws.Range("L131").Value = 13000
$ws.Range("M131").Value = 12500
$ws.Range("N131").Value = '$/caja 50 unidades'
$ws.Range("O131").Value = 'Región de Arica y Parinacota'
$ws.Range("P131").Value = 250
$ws.Range("Q131").Value = 50

# Row 132
$ws.Range("D132").Value = 44817
$ws.Range("H132").Value = 'Huracán'
$ws.Range("J132").Value = 220
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 16000
$ws.Range("M132").Value = 15455
$ws.Range("N132").Value = '$/caja 60 unidades'
$ws.Range("O132").Value = 'Región de Arica y Parinacota'
$ws.Range("P132").Value = 258
$ws.Range("Q132").Value = 60

# Row 133
$ws.Range("D133").Value = 44729
$ws.Range("H133").Value = 'Sin especificar'
$ws.Range("J133").Value = 180
$ws.Range("K133").Value = 13000
$ws.Range("L133").Value = 14000
$ws.Range("M133").Value = 13556
$ws.Range("N133").Value = '$/caja 50 unidades'
$ws.Range("O133").Value = 'Región de Arica y Parinacota'
$ws.Range("P133").Value = 271
$ws.Range("Q133").Value = 50

# Row 134
$ws.Range("D134").Value = 44168
$ws.Range("H134").Value = 'Sin especificar'
$ws.Range("J134").Value = 200
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 7000
$ws.Range("M134").Value = 6500
$ws.Range("N134").Value = '$/caja 50 unidades'
$ws.Range("O134").Value = 'Región de O''Higgins'
$ws.Range("P134").Value = 130
$ws.Range("Q134").Value = 50

# Row 135
$ws.Range("D135").Value = 44455
$ws.Range("H135").Value = 'Sin especificar'
$ws.Range("J135").Value = 100
$ws.Range("K135").Value = 16000
$ws.Range("L135").Value = 17000
$ws.Range("M135").Value = 16500
$ws.Range("N135").Value = '$/caja 50 unidades'
$ws.Range("O135").Value = 'Región de Arica y Parinacota'
$ws.Range("P135").Value = 330
$ws.Range("Q135").Value = 50

# Row 136
$ws.Range("D136").Value = 44719
$ws.Range("H136").Value = 'Sin especificar'
$ws.Range("J136").Value = 100
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 10000
$ws.Range("M136").Value = 9500
$ws.Range("N136").Value = '$/caja 50 unidades'
$ws.Range("O136").Value = 'Región de Arica y Parinacota'
$ws.Range("P136").Value = 190
$ws.Range("Q136").Value = 50

# Row 137
$ws.Range("D137").Value = 44778
$ws.Range("H137").Value = 'Sin especificar'
$ws.Range("J137").Value = 170
$ws.Range("K137").Value = 19000
$ws.Range("L137").Value = 20000
$ws.Range("M137").Value = 19529
$ws.Range("N137").Value = '$/caja 60 unidades'
$ws.Range("O137").Value = 'Región de Arica y Parinacota'
$ws.Range("P137").Value = 325
$ws.Range("Q137").Value = 60

# Row 138
$ws.Range("D138").Value = 44489
$ws.Range("H138").Value = 'Sin especificar'
$ws.Range("J138").Value = 100
$ws.Range("K138").Value = 8000
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 8500
$ws.Range("N138").Value = '$/caja 50 unidades'
$ws.Range("O138").Value = 'Región de Arica y Parinacota'
$ws.Range("P138").Value = 170
$ws.Range("Q138").Value = 50

# Row 139
$ws.Range("D139").Value = 44314
$ws.Range("H139").Value = 'Sin especificar'
$ws.Range("J139").Value = 100
$ws.Range("K139").Value = 8000
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 8500
$ws.Range("N139").Value = '$/caja 60 unidades'
$ws.Range("O139").Value = 'Región del Maule'
$ws.Range("P139").Value = 142
$ws.Range("Q139").Value = 60

# Row 140
$ws.Range("D140").Value = 44656
$ws.Range("H140").Value = 'Sin especificar'
$ws.Range("J140").Value = 220
$ws.Range("K140").Value = 8000
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 8455
$ws.Range("N140").Value = '$/caja 50 unidades'
$ws.Range("O140").Value = 'Región Metropolitana'
$ws.Range("P140").Value = 169
$ws.Range("Q140").Value = 50

# Row 141
$ws.Range("D141").Value = 44894
$ws.Range("H141").Value = 'Sin especificar'
$ws.Range("J141").Value = 100
$ws.Range("K141").Value = 7000
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = 7250
$ws.Range("N141").Value = '$/caja 50 unidades'
$ws.Range("O141").Value = 'Región de O''Higgins'
$ws.Range("P141").Value = 145
$ws.Range("Q141").Value = 50

# Row 142
$ws.Range("D142").Value = 44799
$ws.Range("H142").Value = 'Sin especificar'
$ws.Range("J142").Value = 100
$ws.Range("K142").Value = 22000
$ws.Range("L142").Value = 24000
$ws.Range("M142").Value = 23000
$ws.Range("N142").Value = '$/caja 50 unidades'
$ws.Range("O142").Value = 'Región de Arica y Parinacota'
$ws.Range("P142").Value = 460
$ws.Range("Q142").Value = 50

# Row 143
$ws.Range("D143").Value = 44306
$ws.Range("H143").Value = 'Sin especificar'
$ws.Range("J143").Value = 100
$ws.Range("K143").Value = 8000
$ws.Range("L143").Value = 9000
$ws.Range("M143").Value = 8500
$ws.Range("N143").Value = '$/caja 60 unidades'
$ws.Range("O143").Value = 'Región de O''Higgins'
$ws.Range("P143").Value = 142
$ws.Range("Q143").Value = 60

# Row 144
$ws.Range("D144").Value = 44938
$ws.Range("H144").Value = 'Sin especificar'
$ws.Range("J144").Value = 150
$ws.Range("K144").Value = 8000
$ws.Range("L144").Value = 9000
$ws.Range("M144").Value = 8333
$ws.Range("N144").Value = '$/caja 50 unidades'
$ws.Range("O144").Value = 'Región de O''Higgins'
$ws.Range("P144").Value = 167
$ws.Range("Q144").Value = 50

# Row 145
$ws.Range("D145").Value = 44777
$ws.Range("H145").Value = 'Sin especificar'
$ws.Range("J145").Value = 220
$ws.Range("K145").Value = 17000
$ws.Range("L145").Value = 18000
$ws.Range("M145").Value = 17545
$ws.Range("N145").Value = '$/caja 50 unidades'
$ws.Range("O145").Value = 'Región de Arica y Parinacota'
$ws.Range("P145").Value = 351
$ws.Range("Q145").Value = 50

# Row 146
$ws.Range("D146").Value = 44698
$ws.Range("H146").Value = 'Sin especificar'
$ws.Range("J146").Value = 450
$ws.Range("K146").Value = 17000
$ws.Range("L146").Value = 18000
$ws.Range("M146").Value = 17556
$ws.Range("N146").Value = '$/caja 60 unidades'
$ws.Range("O146").Value = 'Región de Arica y Parinacota'
$ws.Range("P146").Value = 293
$ws.Range("Q146").Value = 60

# Row 147
$ws.Range("D147").Value = 44565
$ws.Range("H147").Value = 'Sin especificar'
$ws.Range("J147").Value = 100
$ws.Range("K147").Value = 7000
$ws.Range("L147").Value = 8000
$ws.Range("M147").Value = 7500
$ws.Range("N147").Value = '$/caja 50 unidades'
$ws.Range("O147").Value = 'Región de O''Higgins'
$ws.Range("P147").Value = 150
$ws.Range("Q147").Value = 50

# Row 148
$ws.Range("D148").Value = 44343
$ws.Range("H148").Value = 'Sin especificar'
$ws.Range("J148").Value = 100
$ws.Range("K148").Value = 9000
$ws.Range("L148").Value = 10000
$ws.Range("M148").Value = 9500
$ws.Range("N148").Value = '$/caja 50 unidades'
$ws.Range("O148").Value = 'Región de Arica y Parinacota'
$ws.Range("P148").Value = 190
$ws.Range("Q148").Value = 50

# Row 149
$ws.Range("D149").Value = 44848
$ws.Range("H149").Value = 'Sin especificar'
$ws.Range("J149").Value = 310
$ws.Range("K149").Value = 17000
$ws.Range("L149").Value = 19000
$ws.Range("M149").Value = 17968
$ws.Range("N149").Value = '$/caja 50 unidades'
$ws.Range("O149").Value = 'Región de O''Higgins'
$ws.Range("P149").Value = 359
$ws.Range("Q149").Value = 50

# Row 150
$ws.Range("D150").Value = 44518
$ws.Range("H150").Value = 'Sin especificar'
$ws.Range("J150").Value = 450
$ws.Range("K150").Value = 6500
$ws.Range("L150").Value = 7000
$ws.Range("M150").Value = 6722
$ws.Range("N150").Value = '$/caja 60 unidades'
$ws.Range("O150").Value = 'Región de O''Higgins'
$ws.Range("P150").Value = 112
$ws.Range("Q150").Value = 60

# Row 151
$ws.Range("D151").Value = 44204
$ws.Range("H151").Value = 'Sin especificar'
$ws.Range("J151").Value = 200
$ws.Range("K151").Value = 6500
$ws.Range("L151").Value = 7000
$ws.Range("M151").Value = 6750
$ws.Range("N151").Value = '$/caja 60 unidades'
$ws.Range("O151").Value = 'Región de O''Higgins'
$ws.Range("P151").Value = 112
$ws.Range("Q151").Value = 60

# Row 152
$ws.Range("D152").Value = 44484
$ws.Range("H152").Value = 'Sin especificar'
$ws.Range("J152").Value = 450
$ws.Range("K152").Value = 12000
$ws.Range("L152").Value = 13000
$ws.Range("M152").Value = 12556
$ws.Range("N152").Value = '$/caja 50 unidades'
$ws.Range("O152").Value = 'Región de O''Higgins'
$ws.Range("P152").Value = 251
$ws.Range("Q152").Value = 50

# Row 153
$ws.Range("D153").Value = 44726
$ws.Range("H153").Value = 'Sin especificar'
$ws.Range("J153").Value = 310
$ws.Range("K153").Value = 10000
$ws.Range("L153").Value = 11000
$ws.Range("M153").Value = 10484
$ws.Range("N153").Value = '$/caja 50 unidades'
$ws.Range("O153").Value = 'Región de Arica y Parinacota'
$ws.Range("P153").Value = 210
$ws.Range("Q153").Value = 50

# Row 154
$ws.Range("D154").Value = 44602
$ws.Range("H154").Value = 'Sin especificar'
$ws.Range("J154").Value = 170
$ws.Range("K154").Value = 7000
$ws.Range("L154").Value = 7500
$ws.Range("M154").Value = 7235
$ws.Range("N154").Value = '$/caja 50 unidades'
$ws.Range("O154").Value = 'Región del Maule'
$ws.Range("P154").Value = 145
$ws.Range("Q154").Value = 50

# Row 155
$ws.Range("D155").Value = 44645
$ws.Range("H155").Value = 'Sin especificar'
$ws.Range("J155").Value = 220
$ws.Range("K155").Value = 11000
$ws.Range("L155").Value = 12000
$ws.Range("M155").Value = 11545
$ws.Range("N155").Value = '$/caja 60 unidades'
$ws.Range("O155").Value = 'Región Metropolitana'
$ws.Range("P155").Value = 192
$ws.Range("Q155").Value = 60

# Row 156
$ws.Range("D156").Value = 44902
$ws.Range("H156").Value = 'Sin especificar'
$ws.Range("J156").Value = 430
$ws.Range("K156").Value = 6500
$ws.Range("L156").Value = 7000
$ws.Range("M156").Value = 6733
$ws.Range("N156").Value = '$/caja 50 unidades'
$ws.Range("O156").Value = 'Región de O''Higgins'
$ws.Range("P156").Value = 135
$ws.Range("Q156").Value = 50

# Row 157
$ws.Range("D157").Value = 44385
$ws.Range("H157").Value = 'Sin especificar'
$ws.Range("J157").Value = 100
$ws.Range("K157").Value = 9000
$ws.Range("L157").Value = 10000
$ws.Range("M157").Value = 9500
$ws.Range("N157").Value = '$/caja 50 unidades'
$ws.Range("O157").Value = 'Región de Arica y Parinacota'
$ws.Range("P157").Value = 190
$ws.Range("Q157").Value = 50

# Row 158
$ws.Range("D158").Value = 44526
$ws.Range("H158").Value = 'Sin especificar'
$ws.Range("J158").Value = 100
$ws.Range("K158").Value = 7000
$ws.Range("L158").Value = 8000
$ws.Range("M158").Value = 7500
$ws.Range("N158").Value = '$/caja 50 unidades'
$ws.Range("O158").Value = 'Región de Arica y Parinacota'
$ws.Range("P158").Value = 150
$ws.Range("Q158").Value = 50

# Row 159
$ws.Range("D159").Value = 44722
$ws.Range("H159").Value = 'Sin especificar'
$ws.Range("J159").Value = 250
$ws.Range("K159").Value = 9000
$ws.Range("L159").Value = 10000
$ws.Range("M159").Value = 9600
$ws.Range("N159").Value = '$/caja 50 unidades'
$ws.Range("O159").Value = 'Región de Arica y Parinacota'
$ws.Range("P159").Value = 192
$ws.Range("Q159").Value = 50

# Row 160
$ws.Range("D160").Value = 44705
$ws.Range("H160").Value = 'Sin especificar'
$ws.Range("J160").Value = 150
$ws.Range("K160").Value = 18000
$ws.Range("L160").Value = 19000
$ws.Range("M160").Value = 18333
$ws.Range("N160").Value = '$/caja 50 unidades'
$ws.Range("O160").Value = 'Región de Arica y Parinacota'
$ws.Range("P160").Value = 367
$ws.Range("Q160").Value = 50

# Row 161
$ws.Range("D161").Value = 44383
$ws.Range("H161").Value = 'Sin especificar'
$ws.Range("J161").Value = 100
$ws.Range("K161").Value = 9000
$ws.Range("L161").Value = 10000
$ws.Range("M161").Value = 9500
$ws.Range("N161").Value = '$/caja 50 unidades'
$ws.Range("O161").Value = 'Región de Arica y Parinacota'
$ws.Range("P161").Value = 190
$ws.Range("Q161").Value = 50

# Row 162
$ws.Range("D162").Value = 44166
$ws.Range("H162").Value = 'Sin especificar'
$ws.Range("J162").Value = 200
$ws.Range("K162").Value = 6000
$ws.Range("L162").Value = 7000
$ws.Range("M162").Value = 6500
$ws.Range("N162").Value = '$/caja 50 unidades'
$ws.Range("O162").Value = 'Región de O''Higgins'
$ws.Range("P162").Value = 130
$ws.Range("Q162").Value = 50

# Row 163
$ws.Range("D163").Value = 44846
$ws.Range("H163").Value = 'Sin especificar'
$ws.Range("J163").Value = 220
$ws.Range("K163").Value = 18000
$ws.Range("L163").Value = 20000
$ws.Range("M163").Value = 19091
$ws.Range("N163").Value = '$/caja 50 unidades'
$ws.Range("O163").Value = 'Región de O''Higgins'
$ws.Range("P163").Value = 382
$ws.Range("Q163").Value = 50

# Row 164
$ws.Range("D164").Value = 44964
$ws.Range("H164").Value = 'Sin especificar'
$ws.Range("J164").Value = 150
$ws.Range("K164").Value = 8500
$ws.Range("L164").Value = 9000
$ws.Range("M164").Value = 8667
$ws.Range("N164").Value = '$/caja 50 unidades'
$ws.Range("O164").Value = 'Región de O''Higgins'
$ws.Range("P164").Value = 173
$ws.Range("Q164").Value = 50

# Row 165
$ws.Range("D165").Value = 44882
$ws.Range("H165").Value = 'Sin especificar'
$ws.Range("J165").Value = 100
$ws.Range("K165").Value = 9000
$ws.Range("L165").Value = 10000
$ws.Range("M165").Value = 9500
$ws.Range("N165").Value = '$/caja 50 unidades'
$ws.Range("O165").Value = 'Región de O''Higgins'
$ws.Range("P165").Value = 190
$ws.Range("Q165").Value = 50

# Row 166
$ws.Range("D166").Value = 44587
$ws.Range("H166").Value = 'Sin especificar'
$ws.Range("J166").Value = 310
$ws.Range("K166").Value = 10000
$ws.Range("L166").Value = 11000
$ws.Range("M166").Value = 10516
$ws.Range("N166").Value = '$/caja 60 unidades'
$ws.Range("O166").Value = 'Región Metropolitana'
$ws.Range("P166").Value = 175
$ws.Range("Q166").Value = 60

# Row 167
$ws.Range("D167").Value = 44215
$ws.Range("H167").Value = 'Sin especificar'
$ws.Range("J167").Value = 100
$ws.Range("K167").Value = 9000
$ws.Range("L167").Value = 10000
$ws.Range("M167").Value = 9500
$ws.Range("N167").Value = '$/caja 60 unidades'
$ws.Range("O167").Value = 'Región de O''Higgins'
$ws.Range("P167").Value = 158
$ws.Range("Q167").Value = 60

# Row 168
$ws.Range("D168").Value = 44782
$ws.Range("H168").Value = 'Sin especificar'
$ws.Range("J168").Value = 100
$ws.Range("K168").Value = 20000
$ws.Range("L168").Value = 22000
$ws.Range("M168").Value = 21000
$ws.Range("N168").Value = '$/caja 50 unidades'
$ws.Range("O168").Value = 'Región de Arica y Parinacota'
$ws.Range("P168").Value = 420
$ws.Range("Q168").Value = 50

# Row 169
$ws.Range("D169").Value = 44754
$ws.Range("H169").Value = 'Sin especificar'
$ws.Range("J169").Value = 100
$ws.Range("K169").Value = 10000
$ws.Range("L169").Value = 12000
$ws.Range("M169").Value = 11000
$ws.Range("N169").Value = '$/caja 50 unidades'
$ws.Range("O169").Value = 'Región de Arica y Parinacota'
$ws.Range("P169").Value = 220
$ws.Range("Q169").Value = 50

# Row 170
$ws.Range("D170").Value = 44329
$ws.Range("H170").Value = 'Sin especificar'
$ws.Range("J170").Value = 900
$ws.Range("K170").Value = 350
$ws.Range("L170").Value = 12000
$ws.Range("M170").Value = 1633
$ws.Range("N170").Value = '$/caja 50 unidades'
$ws.Range("O170").Value = 'Región de Arica y Parinacota'
$ws.Range("P170").Value = 33
$ws.Range("Q170").Value = 50

# Row 171
$ws.Range("D171").Value = 44540
$ws.Range("H171").Value = 'Huracán'
$ws.Range("J171").Value = 190
$ws.Range("K171").Value = 6000
$ws.Range("L171").Value = 6500
$ws.Range("M171").Value = 6263
$ws.Range("N171").Value = '$/caja 60 unidades'
$ws.Range("O171").Value = 'Región del Maule'
$ws.Range("P171").Value = 104
$ws.Range("Q171").Value = 60

# Row 172
$ws.Range("D172").Value = 44904
$ws.Range("H172").Value = 'Sin especificar'
$ws.Range("J172").Value = 450
$ws.Range("K172").Value = 6000
$ws.Range("L172").Value = 6500
$ws.Range("M172").Value = 6278
$ws.Range("N172").Value = '$/caja 50 unidades'
$ws.Range("O172").Value = 'Región Metropolitana'
$ws.Range("P172").Value = 126
$ws.Range("Q172").Value = 50

# Row 173
$ws.Range("D173").Value = 44944
$ws.Range("H173").Value = 'Sin especificar'
$ws.Range("J173").Value = 350
$ws.Range("K173").Value = 9000
$ws.Range("L173").Value = 10000
$ws.Range("M173").Value = 9429
$ws.Range("N173").Value = '$/caja 50 unidades'
$ws.Range("O173").Value = 'Región Metropolitana'
$ws.Range("P173").Value = 189
$ws.Range("Q173").Value = 50

# Row 174
$ws.Range("D174").Value = 44901
$ws.Range("H174").Value = 'Sin especificar'
$ws.Range("J174").Value = 270
$ws.Range("K174").Value = 5500
$ws.Range("L174").Value = 6000
$ws.Range("M174").Value = 5722
$ws.Range("N174").Value = '$/caja 50 unidades'
$ws.Range("O174").Value = 'Región de O''Higgins'
$ws.Range("P174").Value = 114
$ws.Range("Q174").Value = 50

# Row 175
$ws.Range("D175").Value = 44946
$ws.Range("H175").Value = 'Sin especificar'
$ws.Range("J175").Value = 350
$ws.Range("K175").Value = 4000
$ws.Range("L175").Value = 5000
$ws.Range("M175").Value = 4429
$ws.Range("N175").Value = '$/caja 50 unidades'
$ws.Range("O175").Value = 'Región de O''Higgins'
$ws.Range("P175").Value = 89
$ws.Range("Q175").Value = 50

# Row 176
$ws.Range("D176").Value = 44467
$ws.Range("H176").Value = 'Sin especificar'
$ws.Range("J176").Value = 100
$ws.Range("K176").Value = 14000
$ws.Range("L176").Value = 15000
$ws.Range("M176").Value = 14500
$ws.Range("N176").Value = '$/caja 50 unidades'
$ws.Range("O176").Value = 'Región de Arica y Parinacota'
$ws.Range("P176").Value = 290
$ws.Range("Q176").Value = 50

# Row 177
$ws.Range("D177").Value = 44628
$ws.Range("H177").Value = 'Sin especificar'
$ws.Range("J177").Value = 220
$ws.Range("K177").Value = 11000
$ws.Range("L177").Value = 12000
$ws.Range("M177").Value = 11545
$ws.Range("N177").Value = '$/caja 60 unidades'
$ws.Range("O177").Value = 'Región de Arica y Parinacota'
$ws.Range("P177").Value = 192
$ws.Range("Q177").Value = 60

# Row 178
$ws.Range("D178").Value = 44505
$ws.Range("H178").Value = 'Sin especificar'
$ws.Range("J178").Value = 350
$ws.Range("K178").Value = 6500
$ws.Range("L178").Value = 7000
$ws.Range("M178").Value = 6714
$ws.Range("N178").Value = '$/caja 60 unidades'
$ws.Range("O178").Value = 'Región del Maule'
$ws.Range("P178").Value = 112
$ws.Range("Q178").Value = 60

# Row 179
$ws.Range("D179").Value = 44637
$ws.Range("H179").Value = 'Sin especificar'
$ws.Range("J179").Value = 220
$ws.Range("K179").Value = 12000
$ws.Range("L179").Value = 14000
$ws.Range("M179").Value = 12909
$ws.Range("N179").Value = '$/caja 50 unidades'
$ws.Range("O179").Value = 'Región de O''Higgins'
$ws.Range("P179").Value = 258
$ws.Range("Q179").Value = 50

# Row 180
$ws.Range("D180").Value = 44223
$ws.Range("H180").Value = 'Sin especificar'
$ws.Range("J180").Value = 100
$ws.Range("K180").Value = 9000
$ws.Range("L180").Value = 10000
$ws.Range("M180").Value = 9500
$ws.Range("N180").Value = '$/caja 60 unidades'
$ws.Range("O180").Value = 'Región de O''Higgins'
$ws.Range("P180").Value = 158
$ws.Range("Q180").Value = 60

# Row 181
$ws.Range("D181").Value = 44855
$ws.Range("H181").Value = 'Sin especificar'
$ws.Range("J181").Value = 100
$ws.Range("K181").Value = 14000
$ws.Range("L181").Value = 15000
$ws.Range("M181").Value = 14500
$ws.Range("N181").Value = '$/caja 50 unidades'
$ws.Range("O181").Value = 'Región de O''Higgins'
$ws.Range("P181").Value = 290
$ws.Range("Q181").Value = 50

# Row 182
$ws.Range("D182").Value = 44616
$ws.Range("H182").Value = 'Sin especificar'
$ws.Range("J182").Value = 150
$ws.Range("K182").Value = 9000
$ws.Range("L182").Value = 10000
$ws.Range("M182").Value = 9333
$ws.Range("N182").Value = '$/caja 50 unidades'
$ws.Range("O182").Value = 'Región de O''Higgins'
$ws.Range("P182").Value = 187
$ws.Range("Q182").Value = 50

# Row 183
$ws.Range("D183").Value = 44908
$ws.Range("H183").Value = 'Sin especificar'
$ws.Range("J183").Value = 100
$ws.Range("K183").Value = 6000
$ws.Range("L183").Value = 7000
$ws.Range("M183").Value = 6500
$ws.Range("N183").Value = '$/caja 50 unidades'
$ws.Range("O183").Value = 'Región de O''Higgins'
$ws.Range("P183").Value = 130
$ws.Range("Q183").Value = 50

# Row 184
$ws.Range("D184").Value = 44243
$ws.Range("H184").Value = 'Sin especificar'
$ws.Range("J184").Value = 100
$ws.Range("K184").Value = 10000
$ws.Range("L184").Value = 11000
$ws.Range("M184").Value = 10500
$ws.Range("N184").Value = '$/caja 60 unidades'
$ws.Range("O184").Value = 'Región de O''Higgins'
$ws.Range("P184").Value = 175
$ws.Range("Q184").Value = 60

# Row 185
$ws.Range("D185").Value = 44539
$ws.Range("H185").Value = 'Sin especificar'
$ws.Range("J185").Value = 270
$ws.Range("K185").Value = 5000
$ws.Range("L185").Value = 5500
$ws.Range("M185").Value = 5222
$ws.Range("N185").Value = '$/caja 60 unidades'
$ws.Range("O185").Value = 'Región de O''Higgins'
$ws.Range("P185").Value = 87
$ws.Range("Q185").Value = 60

# Row 186
$ws.Range("D186").Value = 44281
$ws.Range("H186").Value = 'Sin especificar'
$ws.Range("J186").Value = 100
$ws.Range("K186").Value = 9000
$ws.Range("L186").Value = 10000
$ws.Range("M186").Value = 9500
$ws.Range("N186").Value = '$/caja 60 unidades'
$ws.Range("O186").Value = 'Región de O''Higgins'
$ws.Range("P186").Value = 158
$ws.Range("Q186").Value = 60

# Row 187
$ws.Range("D187").Value = 44757
$ws.Range("H187").Value = 'Sin especificar'
$ws.Range("J187").Value = 100
$ws.Range("K187").Value = 9000
$ws.Range("L187").Value = 10000
$ws.Range("M187").Value = 9500
$ws.Range("N187").Value = '$/caja 50 unidades'
$ws.Range("O187").Value = 'Región de Arica y Parinacota'
$ws.Range("P187").Value = 190
$ws.Range("Q187").Value = 50

# Row 188
$ws.Range("D188").Value = 44320
$ws.Range("H188").Value = 'Sin especificar'
$ws.Range("J188").Value = 100
$ws.Range("K188").Value = 8000
$ws.Range("L188").Value = 9000
$ws.Range("M188").Value = 8500
$ws.Range("N188").Value = '$/caja 50 unidades'
$ws.Range("O188").Value = 'Región de Arica y Parinacota'
$ws.Range("P188").Value = 170
$ws.Range("Q188").Value = 50

# Row 189
$ws.Range("D189").Value = 44251
$ws.Range("H189").Value = 'Sin especificar'
$ws.Range("J189").Value = 200
$ws.Range("K189").Value = 9000
$ws.Range("L189").Value = 10000
$ws.Range("M189").Value = 9500
$ws.Range("N189").Value = '$/caja 60 unidades'
$ws.Range("O189").Value = 'Región de O''Higgins'
$ws.Range("P189").Value = 158
$ws.Range("Q189").Value = 60

# Row 190
$ws.Range("D190").Value = 44636
$ws.Range("H190").Value = 'Sin especificar'
$ws.Range("J190").Value = 220
$ws.Range("K190").Value = 11000
$ws.Range("L190").Value = 12000
$ws.Range("M190").Value = 11455
$ws.Range("N190").Value = '$/caja 60 unidades'
$ws.Range("O190").Value = 'Provincia de Huasco'
$ws.Range("P190").Value = 191
$ws.Range("Q190").Value = 60

# Row 191
$ws.Range("D191").Value = 44272
$ws.Range("H191").Value = 'Sin especificar'
$ws.Range("J191").Value = 100
$ws.Range("K191").Value = 9000
$ws.Range("L191").Value = 10000
$ws.Range("M191").Value = 9500
$ws.Range("N191").Value = '$/caja 60 unidades'
$ws.Range("O191").Value = 'Región de O''Higgins'
$ws.Range("P191").Value = 158
$ws.Range("Q191").Value = 60

# New row 192 (appended)
$ws.Range("A192").Value = 11
$ws.Range("B192").Value = 'Vega Monumental Concepción'
$ws.Range("C192").Value = 'Bíobío'
$ws.Range("D192").Value = 44889
$ws.Range("E192").Value = 8
$ws.Range("F192").Value = 100112032
$ws.Range("G192").Value = 'Zapallo italiano'
$ws.Range("H192").Value = 'Sin especificar'
$ws.Range("I192").Value = 'Primera'
$ws.Range("J192").Value = 220
$ws.Range("K192").Value = 5500
$ws.Range("L192").Value = 6000
$ws.Range("M192").Value = 5727
$ws.Range("N192").Value = '$/caja 50 unidades'
$ws.Range("O192").Value = 'Región de O''Higgins'
$ws.Range("P192").Value = 115
$ws.Range("Q192").Value = 50
$ws.Range("R192").Value = 'Hortaliza'

# Ensure date formatting/style matches surrounding D column cells (numFmt 165, style index like row191)
$ws.Range("D192").NumberFormat = $ws.Range("D191").NumberFormat